$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for 07-04-2025
$ws.Range("Z1").Value = "07-04-2025 Status"
$ws.Range("AA1").Value = "07-04-2025 Time"

# Copy header style (bold/border/centered) from an existing header cell (Y1) to the new headers
$ws.Range("Y1").Copy()
$ws.Range("Z1:AA1").PasteSpecial(-4122) # xlPasteFormats

# Attendance data for each student row (rows 2-6): Status = "A", Time = "00:00:00"
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 26).Value = "A"         # column Z = 26
    $ws.Cells.Item($r, 27).Value = "00:00:00"  # column AA = 27
}
